# Mise à jour de certains champs de Modules et de Professeurs
#
# Adds a new "Matières enseignés" column (E) to the Feuil1 header row,
# sets the column widths for C/D/E, and moves the active selection to E6
# (mirrors the authored workbook.xml/sheet1.xml/sharedStrings.xml diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1 -> new shared string "Matières enseignés"
$ws.Range("E1").Value = "Matières enseignés"

# Column widths (C, D, E) as set by the author in the Excel UI.
$ws.Columns.Item(3).ColumnWidth = 26.6666666666667
$ws.Columns.Item(4).ColumnWidth = 14.8333333333333
$ws.Columns.Item(5).ColumnWidth = 30.8333333333333

# Active selection ends up on E6 after the edit.
$ws.Range("E6").Select() | Out-Null
